# Add a "Firewall / AADL Component" label textbox on top of the empty
# purple rectangle that represents the Firewall AADL component in the
# diagram repeated across slides 1-4.

$p = $ppt.ActivePresentation

# EMU -> point conversion (1 pt = 12700 EMU), since Shapes.AddTextbox and
# the Left/Top/Width/Height shape properties are expressed in points.
$emuPerPt = 12700

$left   = 3791923 / $emuPerPt
$top    = 2769030 / $emuPerPt
$width  = 1776448 / $emuPerPt
$height = 584775  / $emuPerPt

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    $shp = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)

    # wrap="none" -> the box grows to fit the text instead of wrapping it
    $shp.TextFrame.WordWrap = 0

    $tr = $shp.TextFrame.TextRange
    $tr.Text = "Firewall`rAADL Component"
    $tr.Font.Size = 16
    $tr.Font.Italic = 1

    # <a:spAutoFit/> - shape resizes (height) to fit the text
    $shp.TextFrame.AutoSize = 1

    # <a:noFill/> on the shape
    $shp.Fill.Visible = 0
}
